$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 823234
$ws.Range("B2").Value = "Test SOF 2"
$ws.Range("E2").Value = 88890
$ws.Range("F2").Value = "Test DEA 5"
$ws.Range("H2").Value = "72818;72809;72802;72800"

# Row 3 updates
$ws.Range("A3").Value = 823234
$ws.Range("B3").Value = "Test SOF 2"
$ws.Range("E3").Value = 76542
$ws.Range("F3").Value = "Test DEA 6"

# Row 4 updates
$ws.Range("A4").Value = 765865
$ws.Range("B4").Value = "Test SOF 3"
$ws.Range("E4").Value = 87569
$ws.Range("F4").Value = "Test DEA 7"
$ws.Range("H4").Value = "72815;72816"

# Update the selected cell in the sheet view
$ws.Range("H2").Select()
